$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing data rows (2-19) down to (3-20)
$ws.Rows("2:2").Insert()

# The inserted row inherits formatting from the row above (header); clear it so the
# new data row matches the plain (unstyled) look of the other data rows.
$ws.Range("A2:R2").ClearFormats()

# Restore the date number format used by the rest of column D
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with the weekly data point
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44490
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112045
$ws.Range("G2").Value = "Zapallo"
$ws.Range("H2").Value = "Camote"
$ws.Range("I2").Value = "1a nueva(o)"
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 450
$ws.Range("L2").Value = 480
$ws.Range("M2").Value = 465
$ws.Range("N2").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O2").Value = "Perú"
$ws.Range("P2").Value = 465
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
